# Add 4 more ADCs in the ADS131M0X family to the "Parts" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the new "MFG PN" column first (matches the shared-string ordering
# produced when these parts were typed in column-by-column).
$ws.Range("B13").Value = "AMC131M01DFMR"
$ws.Range("B14").Value = "AMC131M03DFMR"
$ws.Range("B15").Value = "ADS131B02QPWRQ1"
$ws.Range("B16").Value = "ADS131B04QPWRQ1"

# Manufacturer is the same for all of them.
$ws.Range("A13").Value = "Texas Instruments"
$ws.Range("A14").Value = "Texas Instruments"
$ws.Range("A15").Value = "Texas Instruments"
$ws.Range("A16").Value = "Texas Instruments"

# Distributor part numbers and packages.
$ws.Range("C13").Value = "296-AMC131M01DFMRCT-ND"
$ws.Range("D13").Value = "20-SOIC"

$ws.Range("C14").Value = "296-AMC131M03DFMRCT-ND"
$ws.Range("D14").Value = "20-SOIC"

$ws.Range("C15").Value = "296-ADS131B02QPWRQ1CT-ND"
$ws.Range("D15").Value = "20-TSSOP"

$ws.Range("C16").Value = "296-ADS131B04QPWRQ1CT-ND"
$ws.Range("D16").Value = "20-TSSOP"

# Leave the selection on the last entered cell, as in the authored edit.
$ws.Range("C16").Select() | Out-Null
